$wb = $excel.ActiveWorkbook

# Overview sheet: mark the dac72cc2 file's Status as ready for handoff
# (row 3 = dac72cc2-e4b4-4292-bf5d-863e2319ed89.md)
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# zh-cn detail sheet: same file's Status + refreshed handoff timestamp
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-03-09 06:26:05"
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-03-09 06:26:05"

# de-de detail sheet: same file's Status + refreshed handoff timestamp
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-03-09 06:26:15"
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-03-09 06:26:15"
